$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Anchor on the PSO paragraph (stable anchor - inserting paragraphs after it
# does not reseat it, unlike inserting before a paragraph we still need to
# address later).
$pso = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "PARTICLE SWARM OPTIMIZATIONS (PSO)`r") {
        $pso = $p
        break
    }
}

# Make room for the four new paragraphs that go between the PSO paragraph
# and "NEWTON RAPHSON ".
$pso.Range.InsertParagraphAfter()
$pso.Range.InsertParagraphAfter()
$pso.Range.InsertParagraphAfter()
$pso.Range.InsertParagraphAfter()

$overviewPara = $pso.Next(1)
$bodyPara     = $pso.Next(2)
$blankPara    = $pso.Next(3)
$tabOnlyPara  = $pso.Next(4)
$newtonPara   = $pso.Next(5)

$overviewPara.Range.InsertXML("<w:p $wNs><w:r><w:t>Overview:</w:t></w:r></w:p>")

$bodyPara.Range.InsertXML("<w:p $wNs><w:r><w:tab/><w:t>The process of finding optimal values for the specific parameters of a given system to fulfill all design requirements while considering the lowest possible cost is referred to as an optimization.</w:t></w:r></w:p>")

$blankPara.Range.InsertXML("<w:p $wNs></w:p>")

$tabOnlyPara.Range.InsertXML("<w:p $wNs><w:r><w:tab/></w:r></w:p>")

# Give the "NEWTON RAPHSON " paragraph itself a trailing tab run, in
# addition to its original text run.
$newtonPara.Range.InsertXML("<w:p $wNs><w:r><w:t xml:space=""preserve"">NEWTON RAPHSON </w:t></w:r><w:r><w:tab/></w:r></w:p>")

# Remove the old empty paragraph that used to separate "NEWTON RAPHSON "
# from "Algorithm" (it has now been superseded by the new layout above).
$oldBlank = $newtonPara.Next(1)
$oldBlank.Range.Delete()
